$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-12-24 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-25 Thursday", 1)

# Update the table cells (row/col are 1-based Word table indices; some source
# strings repeat across cells with different replacements, so each cell is
# addressed directly via Tables(1).Cell(row, col) rather than a global replace).
$tbl = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; Old = "37÷4=9, 1";   New = "11÷4=2, 3" },
    @{ Row = 1;  Col = 2; Old = "37÷9=4, 1";   New = "89÷9=9, 8" },
    @{ Row = 1;  Col = 3; Old = "72÷8=9, 0";   New = "48÷3=16, 0" },
    @{ Row = 1;  Col = 4; Old = "92÷5=18, 2";  New = "27÷6=4, 3" },
    @{ Row = 1;  Col = 5; Old = "63÷9=7, 0";   New = "90÷6=15, 0" },

    @{ Row = 5;  Col = 1; Old = "27÷2=13, 1";  New = "14÷3=4, 2" },
    @{ Row = 5;  Col = 2; Old = "56÷3=18, 2";  New = "77÷7=11, 0" },
    @{ Row = 5;  Col = 3; Old = "74÷6=12, 2";  New = "45÷6=7, 3" },
    @{ Row = 5;  Col = 4; Old = "10÷6=1, 4";   New = "61÷9=6, 7" },
    @{ Row = 5;  Col = 5; Old = "95÷2=47, 1";  New = "32÷3=10, 2" },

    @{ Row = 9;  Col = 1; Old = "40÷4=10, 0";  New = "88÷8=11, 0" },
    @{ Row = 9;  Col = 2; Old = "35÷6=5, 5";   New = "65÷2=32, 1" },
    @{ Row = 9;  Col = 3; Old = "73÷6=12, 1";  New = "12÷3=4, 0" },
    @{ Row = 9;  Col = 4; Old = "18÷8=2, 2";   New = "51÷7=7, 2" },
    @{ Row = 9;  Col = 5; Old = "79÷3=26, 1";  New = "70÷7=10, 0" },

    @{ Row = 13; Col = 1; Old = "70÷5=14, 0";  New = "21÷9=2, 3" },
    @{ Row = 13; Col = 2; Old = "54÷9=6, 0";   New = "84÷9=9, 3" },
    @{ Row = 13; Col = 3; Old = "50÷2=25, 0";  New = "50÷7=7, 1" },
    @{ Row = 13; Col = 4; Old = "54÷9=6, 0";   New = "27÷2=13, 1" },
    @{ Row = 13; Col = 5; Old = "96÷9=10, 6";  New = "47÷3=15, 2" },

    @{ Row = 17; Col = 1; Old = "93÷7=13, 2";  New = "87÷7=12, 3" },
    @{ Row = 17; Col = 2; Old = "96÷3=32, 0";  New = "11÷8=1, 3" },
    @{ Row = 17; Col = 3; Old = "88÷6=14, 4";  New = "46÷6=7, 4" },
    @{ Row = 17; Col = 4; Old = "46÷5=9, 1";   New = "89÷6=14, 5" },
    @{ Row = 17; Col = 5; Old = "56÷2=28, 0";  New = "52÷3=17, 1" }
)

foreach ($edit in $edits) {
    $cell = $tbl.Cell($edit.Row, $edit.Col)
    $rng = $cell.Range
    # Replace:=1 (wdReplaceOne) so the substitution stays confined to this
    # cell's range instead of touching every matching cell in the table.
    $rng.Find.Execute($edit.Old, $true, $false, $false, $false, $false,
                       $true, 1, $false, $edit.New, 1)
}
